$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "(nan)" placeholders in theta_se (row 4) and lambda_se (row 6) with
# the computed standard errors now available after pickling the bootstrap
# results used in the replication. Values are written column by column so
# that each column's theta_se is followed by its lambda_se, matching the
# order the underlying results table was iterated over.

$ws.Range("B4").Value = "(0.75)"
$ws.Range("B6").Value = "(0.83)"

$ws.Range("C4").Value = "(0.66)"
$ws.Range("C6").Value = "(0.21)"

$ws.Range("D4").Value = "(1.29)"
$ws.Range("D6").Value = "(0.04)"

$ws.Range("E4").Value = "(0.6)"
$ws.Range("E6").Value = "(0.19)"

$ws.Range("F4").Value = "(0.79)"
$ws.Range("F6").Value = "(0.42)"

$ws.Range("G4").Value = "(1.39)"
$ws.Range("G6").Value = "(0.63)"

$ws.Range("H4").Value = "(2.03)"
$ws.Range("H6").Value = "(1.07)"

$ws.Range("I4").Value = "(1.21)"
$ws.Range("I6").Value = "(0.24)"

$ws.Range("J4").Value = "(4.15)"
$ws.Range("J6").Value = "(3.38)"
